# Rename the existing (only) sheet to "firstsheet" and populate it with
# the sales data, then add a second, blank sheet named "Secondsheet".

$wb = $excel.ActiveWorkbook

# --- First sheet: rename + fill with data -------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "firstsheet"

$ws1.Range("A1").Value = "Total_sales"
$ws1.Range("A2").Value = 10000
$ws1.Range("A3").Value = 20000
$ws1.Range("A4").Value = 30000
$ws1.Range("A5").Value = 40000

# --- Second sheet: add new blank sheet, placed after the first one -----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Secondsheet"
